# Applies the "Översikt LEKSAND" update:
#  - bumps the "Förändrad" date (column C) from 2023-10-05 (45204) to
#    2023-10-06 (45205) for every existing data row
#  - updates the aggregate counters (NT / Rödlistade / Alla arter) on row 2
#  - adds "Gränsticka" to the species list (column R) on row 2
#  - appends two new cases (rows 506 and 507)
#  - extends the sheet dimension / row heights accordingly

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reusable formats, copied from the existing rows so the new rows match
# the workbook's existing styles (date format for B/C, wrap-text for R).
$dateFormat = $ws.Cells.Item(2, 3).NumberFormat

# --- 1. Bump "Förändrad" (column C) for all existing data rows (2-505) ---
$ws.Range("C2:C505").Value2 = 45205

# --- 2. Row 2 aggregate counter updates ---
$ws.Cells.Item(2, 10).Value2 = 20   # J2  NT
$ws.Cells.Item(2, 15).Value2 = 23   # O2  Rödlistade
$ws.Cells.Item(2, 17).Value2 = 49   # Q2  Alla arter

# --- 3. Insert "Gränsticka" into the species list on row 2 (column R) ---
$r2 = $ws.Cells.Item(2, 18)
$species = $r2.Value2
$species = $species -replace "Grantaggsvamp`r`n", "Grantaggsvamp`r`nGränsticka`r`n"
$r2.Value2 = $species
$ws.Rows.Item(2).RowHeight = 15

# --- 4. Give row 505 its standard row height (it becomes a regular row) ---
$ws.Rows.Item(505).RowHeight = 15

# --- 5. Append new row 506 ---
$ws.Cells.Item(506, 1).Value2 = "A 47673-2023"
$ws.Cells.Item(506, 2).Value2 = 45203
$ws.Cells.Item(506, 2).NumberFormat = $dateFormat
$ws.Cells.Item(506, 3).Value2 = 45205
$ws.Cells.Item(506, 3).NumberFormat = $dateFormat
$ws.Cells.Item(506, 4).Value2 = "DALARNAS LÄN"
$ws.Cells.Item(506, 5).Value2 = "LEKSAND"
$ws.Cells.Item(506, 7).Value2 = 2
$ws.Cells.Item(506, 8).Value2 = 0
$ws.Cells.Item(506, 9).Value2 = 0
$ws.Cells.Item(506, 10).Value2 = 0
$ws.Cells.Item(506, 11).Value2 = 0
$ws.Cells.Item(506, 12).Value2 = 0
$ws.Cells.Item(506, 13).Value2 = 0
$ws.Cells.Item(506, 14).Value2 = 0
$ws.Cells.Item(506, 15).Value2 = 0
$ws.Cells.Item(506, 16).Value2 = 0
$ws.Cells.Item(506, 17).Value2 = 0
$ws.Cells.Item(506, 18).WrapText = $true
$ws.Rows.Item(506).RowHeight = 15

# --- 6. Append new row 507 (kept as the sheet's new last row) ---
$ws.Cells.Item(507, 1).Value2 = "A 47655-2023"
$ws.Cells.Item(507, 2).Value2 = 45203
$ws.Cells.Item(507, 2).NumberFormat = $dateFormat
$ws.Cells.Item(507, 3).Value2 = 45205
$ws.Cells.Item(507, 3).NumberFormat = $dateFormat
$ws.Cells.Item(507, 4).Value2 = "DALARNAS LÄN"
$ws.Cells.Item(507, 5).Value2 = "LEKSAND"
$ws.Cells.Item(507, 7).Value2 = 2.3
$ws.Cells.Item(507, 8).Value2 = 0
$ws.Cells.Item(507, 9).Value2 = 0
$ws.Cells.Item(507, 10).Value2 = 0
$ws.Cells.Item(507, 11).Value2 = 0
$ws.Cells.Item(507, 12).Value2 = 0
$ws.Cells.Item(507, 13).Value2 = 0
$ws.Cells.Item(507, 14).Value2 = 0
$ws.Cells.Item(507, 15).Value2 = 0
$ws.Cells.Item(507, 16).Value2 = 0
$ws.Cells.Item(507, 17).Value2 = 0
$ws.Cells.Item(507, 18).WrapText = $true

Write-Host "done"
